$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 62
$ws.Range("H62").Value = 2611.5
$ws.Range("I62").Value = 2616.2727
$ws.Range("J62").Value = 2594
$ws.Range("K62").Value = 2616.2727
$ws.Range("L62").Value = 2594
$ws.Range("M62").Value = -1992.2727
$ws.Range("N62").Value = -3842

# Row 64
$ws.Range("H64").Value = 2700
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 2700
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 2700
$ws.Range("N64").Value = -3196
$ws.Range("M64").ClearContents()

# Row 65
$ws.Range("H65").Value = 2611.5
$ws.Range("I65").Value = 2616.2727
$ws.Range("J65").Value = 2594
$ws.Range("K65").Value = 13081.3635
$ws.Range("L65").Value = 12970
$ws.Range("M65").Value = -9961.363499999999
$ws.Range("N65").Value = -19210

# Row 67
$ws.Range("H67").Value = 2700
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 2700
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 2700
$ws.Range("N67").Value = -4416
$ws.Range("M67").ClearContents()

# Row 76
$ws.Range("H76").Value = 46259.477
$ws.Range("I76").Value = 46259.477
$ws.Range("K76").Value = 46259.477
$ws.Range("M76").Value = -45944.477

# Row 79
$ws.Range("H79").Value = 46259.477
$ws.Range("I79").Value = 46259.477
$ws.Range("K79").Value = 46259.477
$ws.Range("M79").Value = -45167.477

# Row 132
$ws.Range("H132").Value = 1906486.6
$ws.Range("I132").Value = 1958679.5
$ws.Range("J132").Value = 1450
$ws.Range("K132").Value = 5876038.5
$ws.Range("L132").Value = 4350
$ws.Range("M132").Value = -5873508.5
$ws.Range("N132").Value = -9410

# Row 135
$ws.Range("H135").Value = 23259678
$ws.Range("I135").Value = 30303714
$ws.Range("J135").Value = 14357.2
$ws.Range("K135").Value = 272733426
$ws.Range("L135").Value = 129214.8
$ws.Range("M135").Value = -272730891
$ws.Range("N135").Value = -134284.8

# Row 138
$ws.Range("H138").Value = 7097.91
$ws.Range("I138").Value = 4245
$ws.Range("J138").Value = 7616.621
$ws.Range("K138").Value = 12735
$ws.Range("L138").Value = 22849.863
$ws.Range("M138").Value = -7595
$ws.Range("N138").Value = -33129.863

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 38257.816
$ws.Range("I32").Value = 37129.65
$ws.Range("K32").Value = 37129.65
$ws.Range("M32").Value = -36842.65

# Row 61
$ws.Range("H61").Value = 2006.3438
$ws.Range("I61").Value = 1528.52
$ws.Range("J61").Value = 3712.8572
$ws.Range("K61").Value = 1528.52
$ws.Range("L61").Value = 3712.8572
$ws.Range("M61").Value = -1316.52
$ws.Range("N61").Value = -4136.8572

# Row 88
$ws.Range("H88").Value = 628416.1
$ws.Range("I88").Value = 1432750.9
$ws.Range("J88").Value = 2822.4443
$ws.Range("K88").Value = 1432750.9
$ws.Range("L88").Value = 2822.4443
$ws.Range("M88").Value = -1432344.9
$ws.Range("N88").Value = -3634.4443

# Row 91
$ws.Range("H91").Value = 628416.1
$ws.Range("I91").Value = 1432750.9
$ws.Range("J91").Value = 2822.4443
$ws.Range("K91").Value = 1432750.9
$ws.Range("L91").Value = 2822.4443
$ws.Range("M91").Value = -1431346.9
$ws.Range("N91").Value = -5630.4443

# Row 136
$ws.Range("H136").Value = 2006.3438
$ws.Range("I136").Value = 1528.52
$ws.Range("J136").Value = 3712.8572
$ws.Range("K136").Value = 4585.559999999999
$ws.Range("L136").Value = 11138.5716
$ws.Range("M136").Value = -2035.559999999999
$ws.Range("N136").Value = -16238.5716

$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value = 2391.8096
$ws.Range("I86").Value = 2936.2727
$ws.Range("J86").Value = 1792.9
$ws.Range("K86").Value = 2936.2727
$ws.Range("L86").Value = 1792.9
$ws.Range("M86").Value = -1813.2727
$ws.Range("N86").Value = -4038.9

# Row 89
$ws.Range("H89").Value = 2391.8096
$ws.Range("I89").Value = 2936.2727
$ws.Range("J89").Value = 1792.9
$ws.Range("K89").Value = 14681.3635
$ws.Range("L89").Value = 8964.5
$ws.Range("M89").Value = -9065.363499999999
$ws.Range("N89").Value = -20196.5

# Row 105
$ws.Range("H105").Value = 3155
$ws.Range("I105").Value = 3155
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 3155
$ws.Range("L105").Value = 0
$ws.Range("M105").Value = -1408
$ws.Range("N105").ClearContents()

# Row 134
$ws.Range("H134").Value = 4287.478
$ws.Range("I134").Value = 4687.1353
$ws.Range("J134").Value = 2644.4443
$ws.Range("K134").Value = 14061.4059
$ws.Range("L134").Value = 7933.3329
$ws.Range("M134").Value = -11526.4059
$ws.Range("N134").Value = -13003.3329

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 25005042
$ws.Range("I31").Value = 5282.303
$ws.Range("J31").Value = 142861060
$ws.Range("K31").Value = 5282.303
$ws.Range("L31").Value = 142861060
$ws.Range("M31").Value = -4987.303
$ws.Range("N31").Value = -142861650

# Row 34
$ws.Range("H34").Value = 25005042
$ws.Range("I34").Value = 5282.303
$ws.Range("J34").Value = 142861060
$ws.Range("K34").Value = 5282.303
$ws.Range("L34").Value = 142861060
$ws.Range("M34").Value = -5080.303
$ws.Range("N34").Value = -142861464

# Row 62
$ws.Range("H62").Value = 55559440
$ws.Range("I62").Value = 4433
$ws.Range("J62").Value = 111114450
$ws.Range("K62").Value = 4433
$ws.Range("L62").Value = 111114450
$ws.Range("M62").Value = -3809
$ws.Range("N62").Value = -111115698

# Row 65
$ws.Range("H65").Value = 55559440
$ws.Range("I65").Value = 4433
$ws.Range("J65").Value = 111114450
$ws.Range("K65").Value = 22165
$ws.Range("L65").Value = 555572250
$ws.Range("M65").Value = -19045
$ws.Range("N65").Value = -555578490

# Row 132
$ws.Range("H132").Value = 4809182.5
$ws.Range("I132").Value = 1464
$ws.Range("J132").Value = 25001600
$ws.Range("K132").Value = 4392
$ws.Range("L132").Value = 75004800
$ws.Range("M132").Value = -1862
$ws.Range("N132").Value = -75009860

# Row 134
$ws.Range("H134").Value = 3124.1177
$ws.Range("I134").Value = 3501.5386
$ws.Range("J134").Value = 1897.5
$ws.Range("K134").Value = 10504.6158
$ws.Range("L134").Value = 5692.5
$ws.Range("M134").Value = -7969.6158
$ws.Range("N134").Value = -10762.5

$ws = $wb.Worksheets.Item("CUL")
# Row 121
$ws.Range("H121").Value = 3265883.5
$ws.Range("I121").Value = 33473
$ws.Range("J121").Value = 4478037.5
$ws.Range("K121").Value = 100419
$ws.Range("L121").Value = 13434112.5
$ws.Range("M121").Value = -99109
$ws.Range("N121").Value = -13436732.5

$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 10855275
$ws.Range("I70").Value = 13081067
$ws.Range("J70").Value = 4538.5
$ws.Range("K70").Value = 13081067
$ws.Range("L70").Value = 4538.5
$ws.Range("M70").Value = -13080797
$ws.Range("N70").Value = -5078.5

# Row 73
$ws.Range("H73").Value = 10855275
$ws.Range("I73").Value = 13081067
$ws.Range("J73").Value = 4538.5
$ws.Range("K73").Value = 13081067
$ws.Range("L73").Value = 4538.5
$ws.Range("M73").Value = -13080131
$ws.Range("N73").Value = -6410.5

# Row 80
$ws.Range("H80").Value = 12000
$ws.Range("I80").Value = 9000
$ws.Range("K80").Value = 9000
$ws.Range("M80").Value = -8002

# Row 83
$ws.Range("H83").Value = 12000
$ws.Range("I83").Value = 9000
$ws.Range("K83").Value = 45000
$ws.Range("M83").Value = -40008

$ws = $wb.Worksheets.Item("LTW")
# Row 136
$ws.Range("H136").Value = 5591.3213
$ws.Range("I136").Value = 6309.409
$ws.Range("J136").Value = 2958.3333
$ws.Range("K136").Value = 18928.227
$ws.Range("L136").Value = 8874.999899999999
$ws.Range("M136").Value = -16378.227
$ws.Range("N136").Value = -13974.9999
